$d = $word.ActiveDocument

# --- Text fixes in the document body ---

# Fix typo: avverkningsamnälda -> avverkningsanmälda
$d.Content.Find.Execute("avverkningsamnälda", $true, $false, $false, $false, $false, $true, 1, $false, "avverkningsanmälda", 2) | Out-Null

# --- Update the date stamp that lives in the (first-page) header ---
foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            $hdr.Range.Find.Execute("2023-11-03", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-13", 2) | Out-Null
        }
    }
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $ftr.Range.Find.Execute("2023-11-03", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-13", 2) | Out-Null
        }
    }
}

# --- Stamp the Swedish editing language onto every (non-numbering) style ---
# Mirrors what Word does when the document's proofing / editing language is
# changed to Swedish: every paragraph/character/table style's rPr gets an
# explicit <w:lang w:val="sv-SE" .../> (eastAsia/bidi stay at their
# existing en-US / ar-SA values).
foreach ($s in $d.Styles) {
    if ($s.Type -ne 4) {
        $f = $s.Font
        $f.LanguageID = "sv-SE"
        $f.LanguageIDFarEast = "en-US"
        $f.LanguageIDOther = "ar-SA"
    }
}
